# CATB_YR_FIN.xlsx update: add the new FY2018 (period ending 2018-12-31)
# column to each of the three statements (Income Statement, Balance Sheet,
# Cash Flow Statement) by inserting a new column D and shifting the
# existing years one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at D - this shifts the previous D:K data to
#    E:L (matches the OOXML diff: dimension A5:K102 -> A5:L102).
$ws.Columns("D").Insert()

# 2) The freshly inserted column D inherits formatting from column C
#    (General / label font). Re-pull the correct per-row number formats /
#    fonts by copying them across from column E (which now holds what used
#    to be column D), so column D matches column E's look (date header
#    style, "#,##0" value style, etc.) before we fill in the new values.
$ws.Columns("E").Copy() | Out-Null
$ws.Columns("D").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Write the new FY2018 column values. Each "Period Ending" header row
#    gets the new date; each data row gets its new leading value following
#    the same pattern already used throughout the sheet ("NA" placeholder
#    where unknown/not applicable, 0 where the series is all zero, or the
#    new reported figure).

# -- Period Ending headers (row 7: Income Statement, 38: Balance Sheet, 80: Cash Flow Statement)
$dateRows = 7, 38, 80
foreach ($r in $dateRows) {
    $ws.Cells.Item($r, 4).Value = 43465   # 31-Dec-2018 (serial date)
}

# -- "NA" placeholder rows
$naRows = 8, 9, 10, 18, 20, 32
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 4).Value = "NA"
}

# -- Numeric rows: new FY2018 figure
$numValues = @{
    12 = 17000
    13 = 0
    14 = 0
    15 = 0
    17 = 26400
    21 = -25700
    22 = 100
    23 = -25900
    24 = 0
    25 = 0
    26 = -25900
    27 = -25900
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    33 = -25900
    34 = 0
    35 = -25900
    41 = 15300
    42 = 22300
    43 = 0
    44 = 0
    45 = 1300
    46 = 38900
    47 = 0
    48 = 100
    49 = 0
    50 = 0
    51 = 0
    52 = 200
    53 = 0
    54 = 39200
    57 = 1400
    58 = 0
    59 = 2800
    60 = 4200
    61 = 0
    62 = 100
    63 = 0
    64 = 0
    65 = 0
    66 = 4200
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -197300
    73 = 0
    74 = 0
    75 = 0
    76 = 34900
    77 = 0
    81 = -25900
    83 = 100
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = -23500
    91 = 0
    92 = 0
    93 = 0
    94 = -21900
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 44300
    101 = 0
    102 = -1100
}
foreach ($r in $numValues.Keys) {
    $ws.Cells.Item($r, 4).Value = $numValues[$r]
}

$wb.Saved = $false
